# Trade #37 closed at 2026-02-17 21:03:12 - unknown UNKNOWN +0.000%
# Also records a brand-new open trade (#98) that was created right after.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: roll up totals after the trade closes + the new trade
# opens.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.56   # Current Capital
$summary.Range("B4").Value = 0.36      # Total P&L $
$summary.Range("B5").Value = 0.11      # Total P&L %
$summary.Range("B6").Value = 65        # Total Trades
$summary.Range("B7").Value = 31        # Winning Trades
$summary.Range("B9").Value = 47.69     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 5).
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.56
$status.Range("D5").Value = 32
$status.Range("E5").Value = 0.25
$status.Range("F5").Value = 0.5600000000000001
$status.Range("G5").Value = 53.12

# ---------------------------------------------------------------------
# All Trades sheet: close out trade #65 (row 66) and append the new
# trade #98 (row 99).
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G66").Value = 0.9399999999999999
$allTrades.Range("H66").Value = "CLOSED"
$allTrades.Range("I66").Value = 8.045999999999999
$allTrades.Range("J66").Value = 0.07000000000000001
$allTrades.Range("K66").Value = 100.56
$allTrades.Range("L66").Value = "early_exit"
$allTrades.Range("M66").Value = 0.14

$allTrades.Range("A99").Value = 98
$allTrades.Range("B99").NumberFormat = "@"
$allTrades.Range("B99").Value = "2026-02-17"
$allTrades.Range("B99").Style = "Normal"
$allTrades.Range("C99").Value = "21:03:05"
$allTrades.Range("D99").Value = "MarketMaking"
$allTrades.Range("E99").Value = "DOWN"
$allTrades.Range("F99").Value = 0.87
$allTrades.Range("H99").Value = "OPEN"
$allTrades.Range("I99").Value = 0
$allTrades.Range("J99").Value = 0
$allTrades.Range("K99").Value = 100.4910412885904
$allTrades.Range("M99").Value = 0
$allTrades.Range("N99").Value = 0
$allTrades.Range("O99").Value = 0
$allTrades.Range("P99").Value = 0.6
$allTrades.Range("Q99").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet: same two updates, but this sheet's columns after
# "Capital After" are ordered differently (Entry Slippage, Exit
# Slippage, Confidence, Entry Reason, Exit Reason, Duration).
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Range("G33").Value = 0.9399999999999999
$marketMaking.Range("H33").Value = "CLOSED"
$marketMaking.Range("I33").Value = 8.045999999999999
$marketMaking.Range("J33").Value = 0.07000000000000001
$marketMaking.Range("K33").Value = 100.56
$marketMaking.Range("P33").Value = "early_exit"
$marketMaking.Range("Q33").Value = 0.14

$marketMaking.Range("A66").Value = 98
$marketMaking.Range("B66").NumberFormat = "@"
$marketMaking.Range("B66").Value = "2026-02-17"
$marketMaking.Range("B66").Style = "Normal"
$marketMaking.Range("C66").Value = "21:03:05"
$marketMaking.Range("D66").Value = "MarketMaking"
$marketMaking.Range("E66").Value = "DOWN"
$marketMaking.Range("F66").Value = 0.87
$marketMaking.Range("H66").Value = "OPEN"
$marketMaking.Range("I66").Value = 0
$marketMaking.Range("J66").Value = 0
$marketMaking.Range("K66").Value = 100.4910412885904
$marketMaking.Range("L66").Value = 0
$marketMaking.Range("M66").Value = 0
$marketMaking.Range("N66").Value = 0.6
$marketMaking.Range("O66").Value = "Normal spread capture: 19600 bps"
$marketMaking.Range("Q66").Value = 0
